$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. "figures" sheet: just move the selection - select the whole column D
#    (this also drops the stale topLeftCell scroll position and the
#    tabSelected flag once another sheet becomes active).
# ---------------------------------------------------------------------------
$figuresWs = $wb.Worksheets.Item("figures")
$figuresWs.Columns.Item(4).Select()

# ---------------------------------------------------------------------------
# 2. "tables" sheet: add the new "Column1" column (with comments / excluded
#    ids) to the worksheet and its underlying table, then make it the active
#    sheet/tab with C7 selected.
# ---------------------------------------------------------------------------
$tablesWs = $wb.Worksheets.Item("tables")
$tbl = $tablesWs.ListObjects.Item("Table2")
$tbl.Resize($tablesWs.Range("A1:E7"))

$tablesWs.Range("E1").Value = "Column1"
$tablesWs.Range("D2").Value = "H13"
$tablesWs.Range("E2").Value = "Please check also figure of psychological symptoms - is this really the same variable and does it make sense to use a different way of visualising? "
$tablesWs.Range("D3").Value = "HC3, HC4"
$tablesWs.Range("D4").Value = "D6"
$tablesWs.Range("D5").Value = "NOT SURVEY"
$tablesWs.Range("D6").Value = "H6, H12"
$tablesWs.Range("D7").Value = "HC3"

$tablesWs.PageSetup.PaperSize = 9
$tablesWs.PageSetup.Orientation = 1

$tablesWs.Activate()
$tablesWs.Range("C7").Select()

# ---------------------------------------------------------------------------
# 3. Add a new "general" sheet in front of the existing sheets and fill it
#    with the six general review comments.
# ---------------------------------------------------------------------------
$generalWs = $wb.Worksheets.Add()
$generalWs.Name = "general"

$generalWs.Range("A1").Value = "1. isn't it general practice to mention both the absolute and the relative value of responses? In some graphs the y-axis shows the absolute number, but it is not explicitly used in the balk itself. "
$generalWs.Range("A2").Value = "2. In some cases I felt like the title didn't really match the variables presented. E.g. satisfaction with information: always/sometimes/never - but lets wait what KCE finds about these figures to change it"
$generalWs.Range("A3").Value = "3. I don't see how it is possible that the N for the EQ-5D questions are differing between the different levels "
$generalWs.Range("A4").Value = "4. In all figures you say 'number of patients"" while I prefer to use the term 'respondent' can you change it in every figure? "
$generalWs.Range("A5").Value = "5. If there is ""i don't know"" as a answer option, put this one on the bottem, because its not really one of the real categories "
$generalWs.Range("A6").Value = "6. make all 'i don't know' categories grey "

# Put "general" in front of "figures" (its default Add() position lands
# next to whichever sheet is currently active).
$generalWs.Move($figuresWs)

# Re-apply the tables tab as the active/selected sheet & selection, since
# adding/moving the new sheet shifts activation back to itself. Re-fetch
# the sheet reference by name since the sheet collection was reordered.
$tablesWs = $wb.Worksheets.Item("tables")
$tablesWs.Activate()
$tablesWs.Range("C7").Select()
